$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value2 = 29.41996
$ws.Cells.Item(2, 8).Value2 = 88.25988000000001
$ws.Cells.Item(2, 9).Value2 = 0.6207199949605289
$ws.Cells.Item(2, 10).Value2 = 0.66829493802317
$ws.Cells.Item(2, 13).Value2 = 0.1030276666666667
$ws.Cells.Item(2, 14).Value2 = 0.309083
$ws.Cells.Item(2, 15).Value2 = 0.005678588141197309
$ws.Cells.Item(2, 16).Value2 = 0.005869434938871914
$ws.Cells.Item(2, 17).Value2 = 3.031069832226667
$ws.Cells.Item(2, 18).Value2 = 27.27962849004
$ws.Cells.Item(2, 19).Value2 = 0.003524813202386913
$ws.Cells.Item(2, 20).Value2 = 0.003922513658704435

$ws.Cells.Item(3, 7).Value2 = 29.41996
$ws.Cells.Item(3, 8).Value2 = 88.25988000000001
$ws.Cells.Item(3, 9).Value2 = 0.6207199949605289
$ws.Cells.Item(3, 10).Value2 = 0.66829493802317
$ws.Cells.Item(3, 14).Value2 = 48.75522599999999
$ws.Cells.Item(3, 15).Value2 = 0.8957491941808339
$ws.Cells.Item(3, 16).Value2 = 0.9258536604633588
$ws.Cells.Item(3, 17).Value2 = 478.1255995703199
$ws.Cells.Item(3, 18).Value2 = 4303.130396132879
$ws.Cells.Item(3, 19).Value2 = 0.556009435297825
$ws.Cells.Item(3, 20).Value2 = 0.6187433146378855

$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 7).Value2 = 29.41996
$ws.Cells.Item(4, 8).Value2 = 88.25988000000001
$ws.Cells.Item(4, 9).Value2 = 0.6207199949605289
$ws.Cells.Item(4, 10).Value2 = 0.66829493802317
$ws.Cells.Item(4, 11).Value2 = 2
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 1.7697965
$ws.Cells.Item(4, 14).Value2 = 3.539593
$ws.Cells.Item(4, 15).Value2 = 0.09754608390528599
$ws.Cells.Item(4, 16).Value2 = 0.06721628437535049
$ws.Cells.Item(4, 17).Value2 = 52.06734223814001
$ws.Cells.Item(4, 18).Value2 = 312.40405342884
$ws.Cells.Item(4, 19).Value2 = 0.06054880471010845
$ws.Cells.Item(4, 20).Value2 = 0.04492030260077262

$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 7).Value2 = 29.41996
$ws.Cells.Item(5, 8).Value2 = 88.25988000000001
$ws.Cells.Item(5, 9).Value2 = 0.6207199949605289
$ws.Cells.Item(5, 10).Value2 = 0.66829493802317
$ws.Cells.Item(5, 11).Value2 = 1
$ws.Cells.Item(5, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(5, 13).Value2 = 0.01861733333333333
$ws.Cells.Item(5, 14).Value2 = 0.055852
$ws.Cells.Item(5, 15).Value2 = 0.001026133772682911
$ws.Cells.Item(5, 16).Value2 = 0.001060620222418814
$ws.Cells.Item(5, 17).Value2 = 0.5477212019733334
$ws.Cells.Item(5, 18).Value2 = 4.929490817760001
$ws.Cells.Item(5, 19).Value2 = 0.0006369417502085649
$ws.Cells.Item(5, 20).Value2 = 0.0007088071258075018

$ws.Cells.Item(6, 9).Value2 = 0.1515698101047853
$ws.Cells.Item(6, 10).Value2 = 0.1631868437822795
$ws.Cells.Item(6, 13).Value2 = 0.1030276666666667
$ws.Cells.Item(6, 14).Value2 = 0.309083
$ws.Cells.Item(6, 15).Value2 = 0.005678588141197309
$ws.Cells.Item(6, 16).Value2 = 0.005869434938871914
$ws.Cells.Item(6, 17).Value2 = 0.7401383596707778
$ws.Cells.Item(6, 18).Value2 = 6.661245237037001
$ws.Cells.Item(6, 19).Value2 = 0.0008607025262245616
$ws.Cells.Item(6, 20).Value2 = 0.0009578145624599442

$ws.Cells.Item(7, 9).Value2 = 0.1515698101047853
$ws.Cells.Item(7, 10).Value2 = 0.1631868437822795
$ws.Cells.Item(7, 14).Value2 = 48.75522599999999
$ws.Cells.Item(7, 15).Value2 = 0.8957491941808339
$ws.Cells.Item(7, 16).Value2 = 0.9258536604633588
$ws.Cells.Item(7, 19).Value2 = 0.1357685352635034
$ws.Cells.Item(7, 20).Value2 = 0.1510871366552858

$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 9).Value2 = 0.1515698101047853
$ws.Cells.Item(8, 10).Value2 = 0.1631868437822795
$ws.Cells.Item(8, 11).Value2 = 2
$ws.Cells.Item(8, 12).Value2 = 1
$ws.Cells.Item(8, 13).Value2 = 1.7697965
$ws.Cells.Item(8, 14).Value2 = 3.539593
$ws.Cells.Item(8, 15).Value2 = 0.09754608390528599
$ws.Cells.Item(8, 16).Value2 = 0.06721628437535049
$ws.Cells.Item(8, 17).Value2 = 12.71400509048783
$ws.Cells.Item(8, 18).Value2 = 76.284030542927
$ws.Cells.Item(8, 19).Value2 = 0.01478504141398965
$ws.Cells.Item(8, 20).Value2 = 0.01096881329798559

$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 9).Value2 = 0.1515698101047853
$ws.Cells.Item(9, 10).Value2 = 0.1631868437822795
$ws.Cells.Item(9, 11).Value2 = 1
$ws.Cells.Item(9, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(9, 13).Value2 = 0.01861733333333333
$ws.Cells.Item(9, 14).Value2 = 0.055852
$ws.Cells.Item(9, 15).Value2 = 0.001026133772682911
$ws.Cells.Item(9, 16).Value2 = 0.001060620222418814
$ws.Cells.Item(9, 17).Value2 = 0.1337446823808889
$ws.Cells.Item(9, 18).Value2 = 1.203702141428
$ws.Cells.Item(9, 19).Value2 = 0.0001555309010676557
$ws.Cells.Item(9, 20).Value2 = 0.0001730792665481855

$ws.Cells.Item(10, 7).Value2 = 0.3873096666666667
$ws.Cells.Item(10, 8).Value2 = 1.161929
$ws.Cells.Item(10, 9).Value2 = 0.008171692087327698
$ws.Cells.Item(10, 10).Value2 = 0.008798009571759262
$ws.Cells.Item(10, 13).Value2 = 0.1030276666666667
$ws.Cells.Item(10, 14).Value2 = 0.309083
$ws.Cells.Item(10, 15).Value2 = 0.005678588141197309
$ws.Cells.Item(10, 16).Value2 = 0.005869434938871914
$ws.Cells.Item(10, 17).Value2 = 0.03990361123411111
$ws.Cells.Item(10, 18).Value2 = 0.359132501107
$ws.Cells.Item(10, 19).Value2 = 0.00004640367378061496
$ws.Cells.Item(10, 20).Value2 = 0.00005163934477301334

$ws.Cells.Item(11, 7).Value2 = 0.3873096666666667
$ws.Cells.Item(11, 8).Value2 = 1.161929
$ws.Cells.Item(11, 9).Value2 = 0.008171692087327698
$ws.Cells.Item(11, 10).Value2 = 0.008798009571759262
$ws.Cells.Item(11, 14).Value2 = 48.75522599999999
$ws.Cells.Item(11, 15).Value2 = 0.8957491941808339
$ws.Cells.Item(11, 16).Value2 = 0.9258536604633588
$ws.Cells.Item(11, 17).Value2 = 6.294456776772665
$ws.Cells.Item(11, 18).Value2 = 56.65011099095399
$ws.Cells.Item(11, 19).Value2 = 0.007319786602317682
$ws.Cells.Item(11, 20).Value2 = 0.008145669366804982

$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 7).Value2 = 0.3873096666666667
$ws.Cells.Item(12, 8).Value2 = 1.161929
$ws.Cells.Item(12, 9).Value2 = 0.008171692087327698
$ws.Cells.Item(12, 10).Value2 = 0.008798009571759262
$ws.Cells.Item(12, 11).Value2 = 2
$ws.Cells.Item(12, 12).Value2 = 1
$ws.Cells.Item(12, 13).Value2 = 1.7697965
$ws.Cells.Item(12, 14).Value2 = 3.539593
$ws.Cells.Item(12, 15).Value2 = 0.09754608390528599
$ws.Cells.Item(12, 16).Value2 = 0.06721628437535049
$ws.Cells.Item(12, 17).Value2 = 0.6854592924828333
$ws.Cells.Item(12, 18).Value2 = 4.112755754897
$ws.Cells.Item(12, 19).Value2 = 0.0007971165619986293
$ws.Cells.Item(12, 20).Value2 = 0.0005913695133124261

$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 7).Value2 = 0.3873096666666667
$ws.Cells.Item(13, 8).Value2 = 1.161929
$ws.Cells.Item(13, 9).Value2 = 0.008171692087327698
$ws.Cells.Item(13, 10).Value2 = 0.008798009571759262
$ws.Cells.Item(13, 11).Value2 = 1
$ws.Cells.Item(13, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(13, 13).Value2 = 0.01861733333333333
$ws.Cells.Item(13, 14).Value2 = 0.055852
$ws.Cells.Item(13, 15).Value2 = 0.001026133772682911
$ws.Cells.Item(13, 16).Value2 = 0.001060620222418814
$ws.Cells.Item(13, 17).Value2 = 0.007210673167555555
$ws.Cells.Item(13, 18).Value2 = 0.064896058508
$ws.Cells.Item(13, 19).Value2 = 0.000008385249230772661
$ws.Cells.Item(13, 20).Value2 = 0.000009331346868842159

$ws.Cells.Item(14, 7).Value2 = 10.122265
$ws.Cells.Item(14, 8).Value2 = 20.24453
$ws.Cells.Item(14, 9).Value2 = 0.2135656295858028
$ws.Cells.Item(14, 10).Value2 = 0.153289545846405
$ws.Cells.Item(14, 13).Value2 = 0.1030276666666667
$ws.Cells.Item(14, 14).Value2 = 0.309083
$ws.Cells.Item(14, 15).Value2 = 0.005678588141197309
$ws.Cells.Item(14, 16).Value2 = 0.005869434938871914
$ws.Cells.Item(14, 17).Value2 = 1.042873344331667
$ws.Cells.Item(14, 18).Value2 = 6.257240065989999
$ws.Cells.Item(14, 19).Value2 = 0.001212751251533277
$ws.Cells.Item(14, 20).Value2 = 0.0008997230161546976

$ws.Cells.Item(15, 7).Value2 = 10.122265
$ws.Cells.Item(15, 8).Value2 = 20.24453
$ws.Cells.Item(15, 9).Value2 = 0.2135656295858028
$ws.Cells.Item(15, 10).Value2 = 0.153289545846405
$ws.Cells.Item(15, 14).Value2 = 48.75522599999999
$ws.Cells.Item(15, 15).Value2 = 0.8957491941808339
$ws.Cells.Item(15, 16).Value2 = 0.9258536604633588
$ws.Cells.Item(15, 17).Value2 = 164.5044392356299
$ws.Cells.Item(15, 18).Value2 = 987.0266354137798
$ws.Cells.Item(15, 19).Value2 = 0.1913012406062053
$ws.Cells.Item(15, 20).Value2 = 0.1419236871326599

$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 7).Value2 = 10.122265
$ws.Cells.Item(16, 8).Value2 = 20.24453
$ws.Cells.Item(16, 9).Value2 = 0.2135656295858028
$ws.Cells.Item(16, 10).Value2 = 0.153289545846405
$ws.Cells.Item(16, 11).Value2 = 2
$ws.Cells.Item(16, 12).Value2 = 1
$ws.Cells.Item(16, 13).Value2 = 1.7697965
$ws.Cells.Item(16, 14).Value2 = 3.539593
$ws.Cells.Item(16, 15).Value2 = 0.09754608390528599
$ws.Cells.Item(16, 16).Value2 = 0.06721628437535049
$ws.Cells.Item(16, 17).Value2 = 17.9143491690725
$ws.Cells.Item(16, 18).Value2 = 71.65739667628999
$ws.Cells.Item(16, 19).Value2 = 0.02083249082286195
$ws.Cells.Item(16, 20).Value2 = 0.01030355370538029

$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 7).Value2 = 10.122265
$ws.Cells.Item(17, 8).Value2 = 20.24453
$ws.Cells.Item(17, 9).Value2 = 0.2135656295858028
$ws.Cells.Item(17, 10).Value2 = 0.153289545846405
$ws.Cells.Item(17, 11).Value2 = 1
$ws.Cells.Item(17, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(17, 13).Value2 = 0.01861733333333333
$ws.Cells.Item(17, 14).Value2 = 0.055852
$ws.Cells.Item(17, 15).Value2 = 0.001026133772682911
$ws.Cells.Item(17, 16).Value2 = 0.001060620222418814
$ws.Cells.Item(17, 17).Value2 = 0.1884495815933333
$ws.Cells.Item(17, 18).Value2 = 1.13069748956
$ws.Cells.Item(17, 19).Value2 = 0.0002191469052022809
$ws.Cells.Item(17, 20).Value2 = 0.000162581992210093

$ws.Cells.Item(18, 7).Value2 = 0.2830933333333334
$ws.Cells.Item(18, 8).Value2 = 0.84928
$ws.Cells.Item(18, 9).Value2 = 0.005972873261555284
$ws.Cells.Item(18, 10).Value2 = 0.006430662776386256
$ws.Cells.Item(18, 13).Value2 = 0.1030276666666667
$ws.Cells.Item(18, 14).Value2 = 0.309083
$ws.Cells.Item(18, 15).Value2 = 0.005678588141197309
$ws.Cells.Item(18, 16).Value2 = 0.005869434938871914
$ws.Cells.Item(18, 17).Value2 = 0.02916644558222223
$ws.Cells.Item(18, 18).Value2 = 0.26249801024
$ws.Cells.Item(18, 19).Value2 = 0.00003391748727194233
$ws.Cells.Item(18, 20).Value2 = 0.00003774435677982456

$ws.Cells.Item(19, 7).Value2 = 0.2830933333333334
$ws.Cells.Item(19, 8).Value2 = 0.84928
$ws.Cells.Item(19, 9).Value2 = 0.005972873261555284
$ws.Cells.Item(19, 10).Value2 = 0.006430662776386256
$ws.Cells.Item(19, 14).Value2 = 48.75522599999999
$ws.Cells.Item(19, 15).Value2 = 0.8957491941808339
$ws.Cells.Item(19, 16).Value2 = 0.9258536604633588
$ws.Cells.Item(19, 17).Value2 = 4.600759815253332
$ws.Cells.Item(19, 18).Value2 = 41.40683833727999
$ws.Cells.Item(19, 19).Value2 = 0.005350196410982395
$ws.Cells.Item(19, 20).Value2 = 0.00595385267072268

$ws.Cells.Item(20, 4).Value = "MuSCs"
$ws.Cells.Item(20, 7).Value2 = 0.2830933333333334
$ws.Cells.Item(20, 8).Value2 = 0.84928
$ws.Cells.Item(20, 9).Value2 = 0.005972873261555284
$ws.Cells.Item(20, 10).Value2 = 0.006430662776386256
$ws.Cells.Item(20, 11).Value2 = 2
$ws.Cells.Item(20, 12).Value2 = 1
$ws.Cells.Item(20, 13).Value2 = 1.7697965
$ws.Cells.Item(20, 14).Value2 = 3.539593
$ws.Cells.Item(20, 15).Value2 = 0.09754608390528599
$ws.Cells.Item(20, 16).Value2 = 0.06721628437535049
$ws.Cells.Item(20, 17).Value2 = 0.5010175905066667
$ws.Cells.Item(20, 18).Value2 = 3.00610554304
$ws.Cells.Item(20, 19).Value2 = 0.0005826303963273109
$ws.Cells.Item(20, 20).Value2 = 0.0004322452578995595

$ws.Cells.Item(21, 4).Value = "Resolving-Mac"
$ws.Cells.Item(21, 7).Value2 = 0.2830933333333334
$ws.Cells.Item(21, 8).Value2 = 0.84928
$ws.Cells.Item(21, 9).Value2 = 0.005972873261555284
$ws.Cells.Item(21, 10).Value2 = 0.006430662776386256
$ws.Cells.Item(21, 11).Value2 = 1
$ws.Cells.Item(21, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(21, 13).Value2 = 0.01861733333333333
$ws.Cells.Item(21, 14).Value2 = 0.055852
$ws.Cells.Item(21, 15).Value2 = 0.001026133772682911
$ws.Cells.Item(21, 16).Value2 = 0.001060620222418814
$ws.Cells.Item(21, 17).Value2 = 0.005270442951111112
$ws.Cells.Item(21, 18).Value2 = 0.04743398656
$ws.Cells.Item(21, 19).Value2 = 0.000006128966973636605
$ws.Cells.Item(21, 20).Value2 = 0.000006820490984191176
